$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.017766216615811729
$ws.Range("B1").Value = -0.017767065727291566

$ws.Range("A2").Value = -0.037945621418336691
$ws.Range("B2").Value = -0.03785085309071267

$ws.Range("A3").Value = -0.0018129233589202719
$ws.Range("B3").Value = -0.0018112501549648527

$ws.Range("A4").Value = -0.01992064693646025
$ws.Range("B4").Value = -0.030519860656117302

$ws.Range("A5").Value = -0.029991570952098208
$ws.Range("B5").Value = -0.029990885721744263
